$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had a header row plus 4 data rows (rows 2-5), one
# per combination of user/metric. Trim it down to just the header and a
# single filtered data row, per the "Filter file and Social Test Data"
# commit: remove rows 3-5 entirely (shifting nothing else up since they
# are the last rows).
$ws.Rows("3:5").Delete() | Out-Null

# Update the remaining data row (row 2) to the new filtered values.
$ws.Range("A2").Value = "Rohit Menon"
$ws.Range("B2").Value = "Posts"
$ws.Range("C2").Value = "Pending Approval"
$ws.Range("D2").Value = "John Smith"
$ws.Range("E2").Value = "Comments"

# Match the saved selection from the edited workbook.
$ws.Range("A3:E5").Select() | Out-Null
